$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1104.3889
$ws.Range("I107").Value = 1073.9375
$ws.Range("J107").Value = 1348
$ws.Range("K107").Value = 1073.9375
$ws.Range("L107").Value = 1348
$ws.Range("M107").Value = 846.0625
$ws.Range("N107").Value = -5188
$ws.Range("H125").Value = 1773.625
$ws.Range("J125").Value = 2663.3333
$ws.Range("L125").Value = 23969.9997
$ws.Range("N125").Value = -28889.9997
$ws.Range("H134").Value = 36940
$ws.Range("J134").Value = 36940
$ws.Range("L134").Value = 36940
$ws.Range("N134").Value = -47080
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 29413284
$ws.Range("I137").Value = 1544.9474
$ws.Range("J137").Value = 66668156
$ws.Range("K137").Value = 4634.8422
$ws.Range("L137").Value = 200004468
$ws.Range("M137").Value = -2084.8422
$ws.Range("N137").Value = -200009568
$ws.Range("H138").Value = 2420.963
$ws.Range("I138").Value = 1987.7812
$ws.Range("J138").Value = 3051.0454
$ws.Range("K138").Value = 5963.3436
$ws.Range("L138").Value = 9153.136200000001
$ws.Range("M138").Value = -823.3436000000002
$ws.Range("N138").Value = -19433.1362
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 882.3333
$ws.Range("I141").Value = 567.625
$ws.Range("J141").Value = 3400
$ws.Range("K141").Value = 1702.875
$ws.Range("L141").Value = 10200
$ws.Range("M141").Value = 3477.125
$ws.Range("N141").Value = -20560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1132348.2
$ws.Range("I2").Value = 1044.6154
$ws.Range("J2").Value = 2263651.8
$ws.Range("K2").Value = 1044.6154
$ws.Range("L2").Value = 2263651.8
$ws.Range("M2").Value = -931.6153999999999
$ws.Range("N2").Value = -2263877.8
$ws.Range("H61").Value = 2315790.2
$ws.Range("I61").Value = 2778741
$ws.Range("J61").Value = 1036.25
$ws.Range("K61").Value = 2778741
$ws.Range("L61").Value = 1036.25
$ws.Range("M61").Value = -2778529
$ws.Range("N61").Value = -1460.25
$ws.Range("H110").Value = 1463.5
$ws.Range("I110").Value = 1222
$ws.Range("J110").Value = 1825.75
$ws.Range("K110").Value = 1222
$ws.Range("L110").Value = 1825.75
$ws.Range("M110").Value = 823
$ws.Range("N110").Value = -5915.75
$ws.Range("H116").Value = 1132348.2
$ws.Range("I116").Value = 1044.6154
$ws.Range("J116").Value = 2263651.8
$ws.Range("K116").Value = 1044.6154
$ws.Range("L116").Value = 2263651.8
$ws.Range("M116").Value = 1249.3846
$ws.Range("N116").Value = -2268239.8
$ws.Range("H132").Value = 808156.0600000001
$ws.Range("I132").Value = 1151244.6
$ws.Range("J132").Value = 72966.21000000001
$ws.Range("K132").Value = 3453733.8
$ws.Range("L132").Value = 218898.63
$ws.Range("M132").Value = -3451203.8
$ws.Range("N132").Value = -223958.63
$ws.Range("H136").Value = 2315790.2
$ws.Range("I136").Value = 2778741
$ws.Range("J136").Value = 1036.25
$ws.Range("K136").Value = 8336223
$ws.Range("L136").Value = 3108.75
$ws.Range("M136").Value = -8333673
$ws.Range("N136").Value = -8208.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1132348.2
$ws.Range("I3").Value = 1044.6154
$ws.Range("J3").Value = 2263651.8
$ws.Range("K3").Value = 1044.6154
$ws.Range("L3").Value = 2263651.8
$ws.Range("M3").Value = -930.6153999999999
$ws.Range("N3").Value = -2263879.8
$ws.Range("H82").Value = 13016.5625
$ws.Range("I82").Value = 4170.375
$ws.Range("J82").Value = 21862.75
$ws.Range("K82").Value = 4170.375
$ws.Range("L82").Value = 21862.75
$ws.Range("M82").Value = -3787.375
$ws.Range("N82").Value = -22628.75
$ws.Range("H85").Value = 13016.5625
$ws.Range("I85").Value = 4170.375
$ws.Range("J85").Value = 21862.75
$ws.Range("K85").Value = 4170.375
$ws.Range("L85").Value = 21862.75
$ws.Range("M85").Value = -2844.375
$ws.Range("N85").Value = -24514.75
$ws.Range("H134").Value = 8048833.5
$ws.Range("I134").Value = 10060730
$ws.Range("J134").Value = 1245.6
$ws.Range("K134").Value = 30182190
$ws.Range("L134").Value = 3736.8
$ws.Range("M134").Value = -30179655
$ws.Range("N134").Value = -8806.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4772.727
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 6375
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 6375
$ws.Range("M4").Value = -388
$ws.Range("N4").Value = -6599
$ws.Range("H31").Value = 579074.2
$ws.Range("I31").Value = 1437.8235
$ws.Range("J31").Value = 1280489.8
$ws.Range("K31").Value = 1437.8235
$ws.Range("L31").Value = 1280489.8
$ws.Range("M31").Value = -1142.8235
$ws.Range("N31").Value = -1281079.8
$ws.Range("H34").Value = 579074.2
$ws.Range("I34").Value = 1437.8235
$ws.Range("J34").Value = 1280489.8
$ws.Range("K34").Value = 1437.8235
$ws.Range("L34").Value = 1280489.8
$ws.Range("M34").Value = -1235.8235
$ws.Range("N34").Value = -1280893.8
$ws.Range("H50").Value = 10122.272
$ws.Range("J50").Value = 11000.25
$ws.Range("L50").Value = 11000.25
$ws.Range("N50").Value = -12250.25
$ws.Range("H68").Value = 17114.572
$ws.Range("J68").Value = 19960.4
$ws.Range("L68").Value = 19960.4
$ws.Range("N68").Value = -21458.4
$ws.Range("H71").Value = 17114.572
$ws.Range("J71").Value = 19960.4
$ws.Range("L71").Value = 59881.2
$ws.Range("N71").Value = -67369.20000000001
$ws.Range("H74").Value = 16883.666
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 16883.666
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 16883.666
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -18631.666
$ws.Range("H77").Value = 16883.666
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 16883.666
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 50650.99800000001
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -59386.99800000001
$ws.Range("H132").Value = 2343
$ws.Range("I132").Value = 2461.7727
$ws.Range("J132").Value = 2016.375
$ws.Range("K132").Value = 7385.3181
$ws.Range("L132").Value = 6049.125
$ws.Range("M132").Value = -4855.3181
$ws.Range("N132").Value = -11109.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1426.6666
$ws.Range("I4").Value = 210
$ws.Range("J4").Value = 1730.8334
$ws.Range("K4").Value = 630
$ws.Range("L4").Value = 5192.5002
$ws.Range("M4").Value = -518
$ws.Range("N4").Value = -5416.5002
$ws.Range("H113").Value = 541.625
$ws.Range("I113").Value = 546.38464
$ws.Range("J113").Value = 536
$ws.Range("K113").Value = 1639.15392
$ws.Range("L113").Value = 1608
$ws.Range("M113").Value = 530.84608
$ws.Range("N113").Value = -5948
$ws.Range("H131").Value = 3134.2554
$ws.Range("I131").Value = 5928.1816
$ws.Range("J131").Value = 2280.5557
$ws.Range("K131").Value = 17784.5448
$ws.Range("L131").Value = 6841.6671
$ws.Range("M131").Value = -12744.5448
$ws.Range("N131").Value = -16921.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2067.419
$ws.Range("J5").Value = 2067.419
$ws.Range("L5").Value = 2067.419
$ws.Range("N5").Value = -2291.419
$ws.Range("H132").Value = 1363.569
$ws.Range("I132").Value = 1370.2766
$ws.Range("J132").Value = 1334.909
$ws.Range("K132").Value = 4110.8298
$ws.Range("L132").Value = 4004.727
$ws.Range("M132").Value = -1580.8298
$ws.Range("N132").Value = -9064.727000000001
$ws.Range("H136").Value = 14733.667
$ws.Range("J136").Value = 14733.667
$ws.Range("L136").Value = 44201.001
$ws.Range("N136").Value = -49301.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2603.25
$ws.Range("I7").Value = 2408.7273
$ws.Range("J7").Value = 3316.5
$ws.Range("K7").Value = 2408.7273
$ws.Range("L7").Value = 3316.5
$ws.Range("M7").Value = -2296.7273
$ws.Range("N7").Value = -3540.5
$ws.Range("H126").Value = 2603.25
$ws.Range("I126").Value = 2408.7273
$ws.Range("J126").Value = 3316.5
$ws.Range("K126").Value = 7226.1819
$ws.Range("L126").Value = 9949.5
$ws.Range("M126").Value = -4756.1819
$ws.Range("N126").Value = -14889.5
$ws.Range("H136").Value = 4425.6484
$ws.Range("I136").Value = 4688.893
$ws.Range("J136").Value = 3606.6667
$ws.Range("K136").Value = 14066.679
$ws.Range("L136").Value = 10820.0001
$ws.Range("M136").Value = -11516.679
$ws.Range("N136").Value = -15920.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 16274000
$ws.Range("I2").Value = 3365333.2
$ws.Range("K2").Value = 3365333.2
$ws.Range("M2").Value = -3365221.2
$ws.Range("H132").Value = 4123.116
$ws.Range("I132").Value = 4517.4736
$ws.Range("K132").Value = 13552.4208
$ws.Range("M132").Value = -11022.4208
$ws.Range("H136").Value = 5803.7104
$ws.Range("I136").Value = 6215.8857
$ws.Range("J136").Value = 995
$ws.Range("K136").Value = 18647.6571
$ws.Range("L136").Value = 2985
$ws.Range("M136").Value = -16097.6571
$ws.Range("N136").Value = -8085
